# Season-record columns (Wins / Losses / Ties) were missing from the
# original scrape; add them to the end of the team statistics table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header formatting (bold, bordered, centered) from A1
# onto the three new header cells so they match the rest of row 1.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Houston's 2002 season record, repeated on every player row.
$wins = 84
$losses = 78
$ties = 0

$lastRow = 43
for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 30).Value = $wins
    $ws.Cells.Item($row, 31).Value = $losses
    $ws.Cells.Item($row, 32).Value = $ties
}
